# Weekly update: a new price record was logged for "Cilantro" at
# "Vega Monumental Concepción" and inserted as row 116 (pushing every
# subsequent record down by one row, which is why the sheet's used range
# grows from A1:R219 to A1:R220).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 116; Excel shifts rows 116:219 down to 117:220.
$ws.Rows(116).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(116, 1).Value  = 11
$ws.Cells.Item(116, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(116, 3).Value  = "Bíobío"
$ws.Cells.Item(116, 4).Value  = 44790
$ws.Cells.Item(116, 5).Value  = 8
$ws.Cells.Item(116, 6).Value  = 100112040
$ws.Cells.Item(116, 7).Value  = "Cilantro"
$ws.Cells.Item(116, 8).Value  = "Sin especificar"
$ws.Cells.Item(116, 9).Value  = "Primera"
$ws.Cells.Item(116, 10).Value = 180
$ws.Cells.Item(116, 11).Value = 7000
$ws.Cells.Item(116, 12).Value = 7500
$ws.Cells.Item(116, 13).Value = 7278
$ws.Cells.Item(116, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(116, 15).Value = "Región Metropolitana"
$ws.Cells.Item(116, 16).Value = 202
$ws.Cells.Item(116, 17).Value = 36
$ws.Cells.Item(116, 18).Value = "Hortaliza"
